$d = $word.ActiveDocument

# 1) Wordage fix: shorten the "Find the P-value..." bullet so it no longer
#    tells students to sketch the t-distribution with the applet.
$d.Content.Find.Execute(
    "Find the P-value and compare it to the level of significance. Sketch the t-distribution using the t-distribution applet.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Find the P-value and compare it to the level of significance.", 2
)

# 2) Mark the first row of the confidence-scale table as the repeating
#    header row.
$table = $d.Tables.Item(1)
$table.Rows.Item(1).HeadingFormat = $true

# 3) Touch the table's preferred width so it re-serializes as a clean
#    integer percentage value instead of a float.
$table.PreferredWidthType = 2
$table.PreferredWidth = 125
